$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Merge "4 " + "sec" runs into a single "4 sec" run (TextBox 28, inside Group 38) ---
$grp38 = $s.Shapes.Item(2)
$textBox28 = $grp38.GroupItems.Item(3)
$textBox28.TextFrame.TextRange.Text = "ZZPLACEHOLDERZZ"
$textBox28.TextFrame.TextRange.Text = "4 sec"
$textBox28.Height = 225520 / 12700.0

# --- Merge "2 " + "sec" runs into a single "2 sec" run (TextBox 31) ---
$textBox31 = $s.Shapes.Item(11)
$textBox31.TextFrame.TextRange.Text = "ZZPLACEHOLDERZZ"
$textBox31.TextFrame.TextRange.Text = "2 sec"

# --- Move Group 34 (the 25% ROI group) slightly to the right ---
$grp34 = $s.Shapes.Item(7)
$grp34.Left = 6316850 / 12700.0

$p.Save()
